# Adds "Mapeamento por core/socket/nodo" labels (slides 6-8) and an
# italic "Work Load" label (slide 9) next to the title of each
# "Resultados MPI" slide: small flipped, auto-fit text boxes positioned
# just below the slide title.

$p = $ppt.ActivePresentation

function Add-Label {
    param(
        [int]$slideIndex,
        [double]$left,
        [double]$top,
        [double]$width,
        [double]$height,
        [string[]]$texts
    )

    $slide = $p.Slides.Item($slideIndex)
    $tb = $slide.Shapes.AddTextbox(1, 0, 0, 1, 1)
    $tb.Flip(0)

    $tb.Left = $left
    $tb.Top = $top
    $tb.Width = $width
    $tb.Height = $height

    $tb.Fill.Visible = 0
    $tb.TextFrame.WordWrap = -1
    $tb.TextFrame.AutoSize = 1

    $tr = $tb.TextFrame.TextRange
    $tr.Text = $texts[0]
    $tr.LanguageID = "en-GB"
    for ($i = 1; $i -lt $texts.Length; $i++) {
        [void]$tr.InsertAfter($texts[$i])
    }
    $tb.TextFrame.TextRange.LanguageID = "en-GB"

    return $tb
}

# Slide 6 - "Mapeamento por core"
[void](Add-Label 6 635.1780395507812 148.5 197.07174682617188 29.081260681152344 @("Mapeamento", " por core"))

# Slide 7 - "Mapeamento por socket"
[void](Add-Label 7 635.1781616210938 148.5 197.07174682617188 50.892208099365234 @("Mapeamento", " por socket"))

# Slide 8 - "Mapeamento por nodo"
[void](Add-Label 8 635.1780395507812 148.5 197.07174682617188 29.081260681152344 @("Mapeamento", " por ", "nodo"))

# Slide 9 - "Work Load" (italic)
$tb9 = Add-Label 9 635.1780395507812 148.5 197.07174682617188 29.081260681152344 @("Work Load")
$tb9.TextFrame.TextRange.Font.Italic = -1
